$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH212"
$ws.Range("C2").Value = "WERKSTATT OEKONOMIE (WORKSHOP ECONOMY)"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"
